# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that get refreshed each time the
# handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-06 17:47:27"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-06 17:47:21"
$zhcn.Range("K2").Value = "2016-09-06 17:47:39"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-09-06 17:47:47"
